$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Call Date / Due Date values (columns E & F) for rows 2 & 3 ---
# Row 2: Call Date 2022-12-10 -> 2022-01-10, Due Date 2022-12-19 -> 2022-01-19
$ws.Range("E2").Value = (Get-Date -Year 2022 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F2").Value = (Get-Date -Year 2022 -Month 1 -Day 19 -Hour 0 -Minute 0 -Second 0)

# Row 3: Call Date 2022-12-10 -> 2022-06-10, Due Date 2022-12-19 -> 2022-06-19
$ws.Range("E3").Value = (Get-Date -Year 2022 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F3").Value = (Get-Date -Year 2022 -Month 6 -Day 19 -Hour 0 -Minute 0 -Second 0)

# Row 4 (E4/F4) stays as-is (2022-12-10 / 2022-12-19)

# --- Add new FX related columns: J (From Currency), K (To Currency), L (Exchange Rate), M (As Of) ---
$ws.Range("J1").Value = "From Currency"
$ws.Range("K1").Value = "To Currency"
$ws.Range("L1").Value = "Exchange Rate "
$ws.Range("M1").Value = "As Of"

# Give the "As Of" date column the same date number format used by the Call Date column
$ws.Range("M2:M4").NumberFormat = $ws.Range("E2").NumberFormat

# Row 2
$ws.Range("J2").Value = "USD"
$ws.Range("K2").Value = "INR"
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = (Get-Date -Year 2022 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 3
$ws.Range("J3").Value = "USD"
$ws.Range("K3").Value = "INR"
$ws.Range("L3").Value = 81
$ws.Range("M3").Value = (Get-Date -Year 2022 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Row 4
$ws.Range("J4").Value = "USD"
$ws.Range("K4").Value = "INR"
$ws.Range("L4").Value = 82
$ws.Range("M4").Value = (Get-Date -Year 2022 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Widen the new "As Of" column to fit its contents
$ws.Columns.Item(13).ColumnWidth = 9.29

# Move the active selection from B5 to A5
$ws.Range("A5").Select() | Out-Null
